# Daily attendance processing
# Normalizes the "Recorded By" (column G) audit list so that the primary
# recorder name is listed before the secondary "System" entry.
#   "backup@backdoor.com, System"   -> "System, backup@backdoor.com"
#   "System, dnasr281@gmail.com"    -> "dnasr281@gmail.com, System"
# Applied across every data row on the active "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$recordedByCol = 7   # Column G : "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $newVal = $val.Replace("backup@backdoor.com, System", "System, backup@backdoor.com")
        $newVal = $newVal.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")

        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
